$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLS Data Series")

# There is no unemployment data for this series before 2003, so remove the
# 2000-2002 rows (rows 2:4) and let everything below shift up.
$ws.Rows("2:4").Select()
$ws.Rows("2:4").Delete()

Write-Output "Done"
